$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dates to be placed in column B, rows 1..6
$dates = @("22/10/2021", "25/10/2021", "26/10/2021", "27/10/2021", "28/10/2021", "29/10/2021")

# --- Step 1: bank the pre-existing date / quote-prefix style (currently on B1) onto a
#     scratch cell far away, so we can re-apply it later without Excel trying to
#     "smart reinterpret" the new (non-date) text we are about to type into B1. ---
$ws.Range("B1").Copy()
$ws.Range("ZZ1").PasteSpecial(-4122)

# --- Step 2: insert a new column before the old column E. This shifts the old
#     column E (with its wrap-text style and its 31.140625 width) one column to
#     the right, turning it into column F - exactly matching the target layout -
#     while leaving a blank new column E in its place. ---
$ws.Columns("E").Insert()

# --- Step 3: clear the cells that must be fully repopulated, back to a clean slate ---
$ws.Range("A1:E6").Clear()
$ws.Range("B1").Clear()

# --- Step 4: make sure C:E start from the same "vertical center" baseline the rest
#     of the sheet uses, before typing the quote-prefixed text values into them ---
$ws.Range("C1:E6").VerticalAlignment = -4108

# --- Step 5: write the new values (in the same order the original author likely
#     typed them, so the shared-string table order is preserved: 09, 20, 00,
#     "Entrega de actividades", "Entrega de conocimientos...", then the six dates) ---
for ($r = 1; $r -le 6; $r++) {
    $ws.Range("C$r").Value = "'09"
    $ws.Range("D$r").Value = "'20"
    $ws.Range("E$r").Value = "'00"
    $ws.Range("A$r").Value = "Entrega de actividades"
    $ws.Range("F$r").Value = "Entrega de conocimientos tras salida del equipo de System Test"
    $ws.Range("B$r").Value = $dates[$r - 1]
}

# --- Step 6: (re)apply formatting ---
# Wrap-text / vertical-center on column A (column F already inherited it from the
# old column E via the Insert in step 2)
$ws.Range("A1:A6").WrapText = $true

# Date-format / quote-prefix style (banked in step 1) onto column B
$ws.Range("ZZ1").Copy()
$ws.Range("B1:B6").PasteSpecial(-4122)

# --- Step 7: remove the scratch cell used to bank the style ---
$ws.Range("ZZ1").Clear()

# --- Step 8: row heights for the six data rows ---
$ws.Rows("1:6").RowHeight = 30

# --- Step 9: selection matches the diff (active cell B6) ---
$ws.Range("B6").Select()
